# Update crypto price/volume data per latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.844.13"
$ws.Range("E2").Value = "  -1.72%  "
$ws.Range("D3").Value = "2.219.37"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'242.68"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "'0.625"
$ws.Range("E6").Value = "  -1.11%  "
$ws.Range("D7").Value = "'72.92"
$ws.Range("E7").Value = "  -5.30%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("D9").Value = "'0.596"
$ws.Range("E9").Value = "  -4.75%  "
$ws.Range("D10").Value = "'42.13"
$ws.Range("E10").Value = "  -6.31%  "
$ws.Range("D11").Value = "'0.0951"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.103"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.89"
$ws.Range("E13").Value = "  -5.64%  "
$ws.Range("D14").Value = "2.556.38"
$ws.Range("E14").Value = "  -1.64%  "
$ws.Range("D15").Value = "'14.28"
$ws.Range("E15").Value = "  -2.77%  "
$ws.Range("D16").Value = "'0.834"
$ws.Range("E16").Value = "  -3.20%  "
$ws.Range("D17").Value = "2.212.15"
$ws.Range("E17").Value = "  -1.80%  "
$ws.Range("D18").Value = "41.729.88"
$ws.Range("E18").Value = "  -1.53%  "
$ws.Range("D19").Value = "'0.0000105"
$ws.Range("E19").Value = "  +3.63%  "
$ws.Range("D20").Value = "'6.20"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'72.72"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'11.17"
$ws.Range("E22").Value = "  +17.84%  "
$ws.Range("D23").Value = "'229.78"
$ws.Range("E23").Value = "  -1.15%  "
$ws.Range("D24").Value = "'2.07"
$ws.Range("E24").Value = "  -8.09%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'11.31"
$ws.Range("E26").Value = "  -2.50%  "
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "'2.19"
$ws.Range("E29").Value = "  -0.84%  "
$ws.Range("D30").Value = "'166.64"
$ws.Range("E30").Value = "  -0.26%  "
$ws.Range("D31").Value = "'20.47"
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'5.60"
$ws.Range("E32").Value = "  +4.28%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0796"
$ws.Range("E33").Value = "  -3.59%  "
$ws.Range("D34").Value = "'30.14"
$ws.Range("E34").Value = "  -4.82%  "
$ws.Range("E35").Value = "  -0.10%  "
$ws.Range("D36").Value = "'0.107"
$ws.Range("E36").Value = "  -10.77%  "
$ws.Range("D37").Value = "'4.27"
$ws.Range("E37").Value = "  -6.77%  "
$ws.Range("D38").Value = "'0.0302"
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("D39").Value = "'13.35"
$ws.Range("E39").Value = "  -8.96%  "
$ws.Range("D40").Value = "'2.12"
$ws.Range("E40").Value = "  -2.92%  "
$ws.Range("D41").Value = "'64.82"
$ws.Range("E41").Value = "  +0.30%  "
$ws.Range("D42").Value = "'5.64"
$ws.Range("E42").Value = "  -2.93%  "
$ws.Range("D43").Value = "'0.197"
$ws.Range("E43").Value = "  -2.75%  "
$ws.Range("D44").Value = "'8.75"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("D45").Value = "'103.63"
$ws.Range("E45").Value = "  -3.91%  "
$ws.Range("D46").Value = "'0.100"
$ws.Range("E46").Value = "  -3.23%  "
$ws.Range("D47").Value = "'2.33"
$ws.Range("E47").Value = "  -2.43%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").Value = "'1.17"
$ws.Range("E48").Value = "  -2.00%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.11"
$ws.Range("E49").Value = "  -2.66%  "
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("D51").Value = "2.427.90"
$ws.Range("E51").Value = "  -1.67%  "

# Reset number formatting on the Price column so forced-text values
# (e.g. "1.00") do not retain a quote-prefix style
$ws.Range("D2:D51").Style = "Normal"
